# "Generate Report for Handoff"
# The localization report is regenerated: the row ordering for the two
# tracked files (3b55afe9... and 743a8dc0...) swaps in every sheet, and the
# 3b55afe9 file's status moves from "Handed back: in sync with en-US" to
# "Ready for handoff" with fresh timestamps / a stale-version error detail.

$wb = $excel.ActiveWorkbook

$idA = "3b55afe9-fca1-4f9d-b841-af26fca2fc20.md"
$idB = "743a8dc0-4c9d-4dbf-9795-7434f645abdb.md"

$pathA = "e2e\3b55afe9-fca1-4f9d-b841-af26fca2fc20.md"
$pathB = "e2e\743a8dc0-4c9d-4dbf-9795-7434f645abdb.md"

# ---------------------------------------------------------------------
# Sheet "Overview": row 2 <-> row 3 now reference the opposite files, and
# the (now) row-3 (3b55afe9) entry is "Ready for handoff".
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = $idB
$wsOverview.Range("A3").Value = $idA

$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-09-06 09:03:46"

foreach ($h in $wsOverview.Hyperlinks) {
    $r = $h.Range.Row
    if ($r -eq 2) {
        $h.TextToDisplay = $pathB
    } elseif ($r -eq 3) {
        $h.TextToDisplay = $pathA
    }
}

# ---------------------------------------------------------------------
# Sheet "zh-cn": same row swap; row 3 (3b55afe9) becomes "Ready for
# handoff" with a new handback datetime and a stale-version error detail;
# column P (Error Detail) widens to fit the message.
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("A2").Value = $idB
$wsZh.Range("G2").Value = "743a8dc0-4c9d-4dbf-9795-7434f645abdb.d548fde11194a939a491cb9b1bc79ede6d3f6d11.zh-cn.xlf"
$wsZh.Range("I2").Value = $idB
$wsZh.Range("J2").Value = "743a8dc0-4c9d-4dbf-9795-7434f645abdb.d548fde11194a939a491cb9b1bc79ede6d3f6d11.zh-cn.xlf"

$wsZh.Range("A3").Value = $idA
$wsZh.Range("C3").Value = "Ready for handoff"
$wsZh.Range("G3").Value = "3b55afe9-fca1-4f9d-b841-af26fca2fc20.783808eb53d60caf257a1f5621c198600b8a227c.zh-cn.xlf"
$wsZh.Range("H3").Value = "2016-09-06 09:03:36"
$wsZh.Range("I3").Value = $idA
$wsZh.Range("J3").Value = "3b55afe9-fca1-4f9d-b841-af26fca2fc20.783808eb53d60caf257a1f5621c198600b8a227c.zh-cn.xlf"
$wsZh.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/677d2ed9ce93812cf171508034449025ff985881/e2e/3b55afe9-fca1-4f9d-b841-af26fca2fc20.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ebef4bf3d5bfc9bd75ac7a54368758ba43ea610e/e2e/3b55afe9-fca1-4f9d-b841-af26fca2fc20.md."

foreach ($h in $wsZh.Hyperlinks) {
    $r = $h.Range.Row
    if ($r -eq 2) {
        $h.TextToDisplay = $idB
    } elseif ($r -eq 3) {
        $h.TextToDisplay = $idA
    }
}

$wsZh.Columns.Item(16).ColumnWidth = 39.17

# ---------------------------------------------------------------------
# Sheet "de-de": identical pattern to "zh-cn".
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("A2").Value = $idB
$wsDe.Range("G2").Value = "743a8dc0-4c9d-4dbf-9795-7434f645abdb.d548fde11194a939a491cb9b1bc79ede6d3f6d11.de-de.xlf"
$wsDe.Range("I2").Value = $idB
$wsDe.Range("J2").Value = "743a8dc0-4c9d-4dbf-9795-7434f645abdb.d548fde11194a939a491cb9b1bc79ede6d3f6d11.de-de.xlf"

$wsDe.Range("A3").Value = $idA
$wsDe.Range("C3").Value = "Ready for handoff"
$wsDe.Range("G3").Value = "3b55afe9-fca1-4f9d-b841-af26fca2fc20.783808eb53d60caf257a1f5621c198600b8a227c.de-de.xlf"
$wsDe.Range("H3").Value = "2016-09-06 09:03:46"
$wsDe.Range("I3").Value = $idA
$wsDe.Range("J3").Value = "3b55afe9-fca1-4f9d-b841-af26fca2fc20.783808eb53d60caf257a1f5621c198600b8a227c.de-de.xlf"
$wsDe.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/677d2ed9ce93812cf171508034449025ff985881/e2e/3b55afe9-fca1-4f9d-b841-af26fca2fc20.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ebef4bf3d5bfc9bd75ac7a54368758ba43ea610e/e2e/3b55afe9-fca1-4f9d-b841-af26fca2fc20.md."

foreach ($h in $wsDe.Hyperlinks) {
    $r = $h.Range.Row
    if ($r -eq 2) {
        $h.TextToDisplay = $idB
    } elseif ($r -eq 3) {
        $h.TextToDisplay = $idA
    }
}

$wsDe.Columns.Item(16).ColumnWidth = 39.17
